$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume updates (GitHub Actions symbol-list refresh).
# Cells are plain text (inline strings) in the source sheet, e.g. "0.1050" or
# "7.30%". Force text format before assignment so exact literal formatting
# (trailing zeros, no scientific notation, literal "%" suffix) is preserved
# instead of Excel auto-coercing the text into a General/Percent number.
$cellUpdates = @{
    "D2" = "321.02"
    "E2" = "7.30%"
    "D3" = "48.78"
    "E3" = "15.72%"
    "D4" = "5.251"
    "E4" = "4.58%"
    "D5" = "0.08091"
    "E5" = "6.87%"
    "E6" = "4.85%"
    "D7" = "1.645"
    "E7" = "2.94%"
    "D8" = "1.208"
    "E8" = "27.83%"
    "D9" = "0.1293"
    "E9" = "8.13%"
    "D10" = "0.1943"
    "E10" = "6.54%"
    "D11" = "0.09447"
    "E11" = "3.68%"
    "D12" = "0.04617"
    "E12" = "9.82%"
    "D13" = "0.1050"
    "E13" = "0.19%"
    "D14" = "0.001339"
    "E14" = "3.77%"
    "D15" = "0.04166"
    "E15" = "1.48%"
    "D16" = "0.005863"
    "E16" = "0.42%"
    "D18" = "2.430"
    "E18" = "1.93%"
    "D19" = "0.3407"
    "E19" = "2.19%"
    "D20" = "8.078"
    "E20" = "-3.71%"
    "D21" = "0.1382"
    "E21" = "-1.90%"
    "D22" = "0.3126"
    "E24" = "8.90%"
    "D26" = "0.0003542"
    "E26" = "-4.88%"
    "D38" = "0.02725"
    "E38" = "12.93%"
    "D39" = "0.05725"
    "E39" = "8.78%"
    "D40" = "0.006305"
    "E40" = "-3.37%"
    "D41" = "0.007866"
    "E41" = "1.89%"
    "D42" = "0.1441"
    "E42" = "8.15%"
    "D43" = "0.007706"
    "E43" = "4.32%"
    "E44" = "3.80%"
    "E45" = "6.52%"
    "D46" = "0.00006887"
    "E46" = "7.87%"
    "D47" = "0.00000000751"
    "E47" = "0.12%"
    "D48" = "0.06192"
    "E48" = "35.39%"
    "D49" = "0.004003"
    "E49" = "-4.70%"
    "D50" = "0.00002102"
    "E50" = "0.12%"
    "D51" = "0.0002002"
    "E51" = "0.12%"
}

foreach ($addr in $cellUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $cellUpdates[$addr]
}
